$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record needs to be inserted at row 64, pushing the
# existing rows 64-81 down to 65-82 (dimension grows from R81 to R82).
$ws.Rows.Item(64).Insert()

# Populate the newly inserted row 64 with the new record's data.
$ws.Range("A64").Value = 1
$ws.Range("B64").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C64").Value = "Arica y Parinacota"
$ws.Range("D64").Value = 44798
$ws.Range("E64").Value = 15
$ws.Range("F64").Value = 100114001
$ws.Range("G64").Value = "Papa"
$ws.Range("H64").Value = "Asterix"
$ws.Range("I64").Value = "1a (guarda)"
$ws.Range("J64").Value = 1000
$ws.Range("K64").Value = 12000
$ws.Range("L64").Value = 13000
$ws.Range("M64").Value = 12500
$ws.Range("N64").Value = '$/saco 25 kilos'
$ws.Range("O64").Value = "Región de Los Lagos"
$ws.Range("P64").Value = 500
$ws.Range("Q64").Value = 25
$ws.Range("R64").Value = "Hortaliza"
